# Generate Report for handback
#
# The "d690e023-de58-...md" file has now been handed back (in sync with
# en-US), so update its status from "Ready for handoff" to
# "Handed back: in sync with en-US" on the Overview sheet and on each
# locale sheet, and record the new handback timestamps on the locale
# sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the d690e023... file ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 3 is the d690e023... file ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "2016-01-11 05:03:27"

# --- de-de sheet: row 3 is the d690e023... file ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "2016-01-11 05:03:54"
